$wb = $excel.ActiveWorkbook

# Mapping of sheet -> row -> column letter -> new value (derived from the target diff).
# A value of $null means the cell should be cleared (the <c> element removed entirely),
# matching rows where a trailing column was dropped from the sheet in the diff.
$updates = @{
    "ALC" = @{
        32 = @{ "H" = 7960; "J" = 7960; "L" = 7960; "N" = -8612 }
        43 = @{ "H" = 3049.75; "I" = 3066.6667; "K" = 3066.6667; "M" = -2997.6667 }
        74 = @{ "H" = 5498.5557; "I" = 4800; "K" = 4800; "M" = -3864 }
        77 = @{ "H" = 5498.5557; "I" = 4800; "K" = 24000; "M" = -19320 }
        113 = @{ "H" = 5879.6; "J" = 6879.6; "L" = 6879.6; "N" = -13387.6 }
        131 = @{ "H" = 3188.5; "I" = 2765.2222; "K" = 8295.6666; "M" = -3255.6666 }
        137 = @{ "H" = 11534.375; "I" = 11533.969; "J" = 11536; "K" = 34601.907; "L" = 34608; "M" = -32051.907; "N" = -39708 }
        138 = @{ "H" = 25004376; "I" = 886.7826; "J" = 58832624; "K" = 2660.3478; "L" = 176497872; "M" = 2479.6522; "N" = -176508152 }
        141 = @{ "H" = 1658.4286; "I" = 1658.4286; "J" = 0; "K" = 4975.2858; "L" = 0; "M" = 204.7142000000003; "N" = $null }
    }
    "ARM" = @{
        2 = @{ "H" = 3263.353; "I" = 2859; "K" = 2859; "M" = -2746 }
        45 = @{ "H" = 3200.4443; "I" = 2358.5; "K" = 2358.5; "M" = -1981.5 }
        74 = @{ "H" = 394373; "I" = 394373; "K" = 394373; "M" = -393499 }
        77 = @{ "H" = 394373; "I" = 394373; "K" = 1971865; "M" = -1967497 }
        97 = @{ "H" = 1512.238; "I" = 1225.5; "K" = 1225.5; "M" = -729.5 }
        110 = @{ "H" = 21303; "I" = 21303; "K" = 21303; "M" = -19258 }
        116 = @{ "H" = 3263.353; "I" = 2859; "K" = 2859; "M" = -565 }
        122 = @{ "H" = 2336.75; "I" = 1641.3684; "K" = 4924.1052; "M" = -2474.1052 }
        130 = @{ "H" = 67464.5; "J" = 67464.5; "L" = 67464.5; "N" = -77504.5 }
    }
    "BSM" = @{
        3 = @{ "H" = 3263.353; "I" = 2859; "K" = 2859; "M" = -2745 }
        7 = @{ "H" = 3250; "J" = 3250; "L" = 3250; "N" = -3476 }
        105 = @{ "H" = 8640; "I" = 16417; "K" = 16417; "M" = -14670 }
        132 = @{ "H" = 139779; "J" = 139779; "L" = 139779; "N" = -149899 }
        134 = @{ "H" = 3600.111; "I" = 893.6667; "K" = 2681.0001; "M" = -146.0001000000002 }
    }
    "CRP" = @{
        31 = @{ "H" = 4457.1333; "I" = 2955.6; "J" = 5207.9; "K" = 2955.6; "L" = 5207.9; "M" = -2660.6; "N" = -5797.9 }
        34 = @{ "H" = 4457.1333; "I" = 2955.6; "J" = 5207.9; "K" = 2955.6; "L" = 5207.9; "M" = -2753.6; "N" = -5611.9 }
        99 = @{ "H" = 4127.095; "I" = 4398.5386; "J" = 3686; "K" = 4398.5386; "L" = 3686; "M" = -2900.5386; "N" = -6682 }
        107 = @{ "H" = 449.8889; "I" = 423.57144; "K" = 423.57144; "M" = 1496.42856 }
        126 = @{ "H" = 4127.095; "I" = 4398.5386; "J" = 3686; "K" = 13195.6158; "L" = 11058; "M" = -10725.6158; "N" = -15998 }
    }
    "CUL" = @{
        34 = @{ "H" = 806.9231; "I" = 151.25; "J" = 1856; "K" = 453.75; "L" = 5568; "M" = -369.75; "N" = -5736 }
        55 = @{ "H" = 12611244; "I" = 2250150; "J" = 33333432; "K" = 6750450; "L" = 100000296; "M" = -6750273; "N" = -100000650 }
        109 = @{ "H" = 8420; "J" = 2123.5; "L" = 6370.5; "N" = -8450.5 }
        131 = @{ "H" = 1479.3334; "I" = 1197.5; "J" = 1704.8; "K" = 3592.5; "L" = 5114.4; "M" = 1447.5; "N" = -15194.4 }
        132 = @{ "H" = 2916.318; "J" = 3718.2; "L" = 33463.8; "N" = -38523.8 }
        134 = @{ "H" = 1021.7895; "I" = 1021.7895; "J" = 0; "K" = 3065.3685; "L" = 0; "M" = 2004.6315; "N" = $null }
        137 = @{ "H" = 2454.1667; "I" = 2339; "J" = 3030; "K" = 7017; "L" = 9090; "M" = -1917; "N" = -19290 }
    }
    "GSM" = @{
        80 = @{ "H" = 2915.6667; "I" = 2798.3333; "J" = 2954.7778; "K" = 2798.3333; "L" = 2954.7778; "M" = -1800.3333; "N" = -4950.7778 }
        83 = @{ "H" = 2915.6667; "I" = 2798.3333; "J" = 2954.7778; "K" = 13991.6665; "L" = 14773.889; "M" = -8999.666499999999; "N" = -24757.889 }
        102 = @{ "H" = 5031.727; "I" = 2586.5; "K" = 2586.5; "M" = -964.5 }
        113 = @{ "H" = 3982.8333; "I" = 2000; "K" = 2000; "M" = 170 }
        122 = @{ "H" = 946.4375; "I" = 799.5484; "K" = 2398.6452; "M" = 51.35480000000007 }
    }
    "LTW" = @{
        7 = @{ "H" = 5890.1665; "I" = 5095.5713; "K" = 5095.5713; "M" = -4983.5713 }
        40 = @{ "H" = 1914.3334; "J" = 5000; "L" = 5000; "N" = -5272 }
        126 = @{ "H" = 5890.1665; "I" = 5095.5713; "K" = 15286.7139; "M" = -12816.7139 }
        136 = @{ "H" = 3373.36; "I" = 2405.8262; "K" = 7217.4786; "M" = -4667.4786 }
    }
    "WVR" = @{
        119 = @{ "H" = 83925; "J" = 83925; "L" = 83925; "N" = -93601 }
        122 = @{ "H" = 42689.348; "I" = 65599.06; "K" = 196797.18; "M" = -194347.18 }
        136 = @{ "H" = 5747.6; "I" = 3996.5; "K" = 11989.5; "M" = -9439.5 }
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $sheetRows = $updates[$sheetName]
    foreach ($rowNum in $sheetRows.Keys) {
        $rowCols = $sheetRows[$rowNum]
        foreach ($colLetter in $rowCols.Keys) {
            $newValue = $rowCols[$colLetter]
            $cellRef = "$colLetter$rowNum"
            if ($null -eq $newValue) {
                $ws.Range($cellRef).ClearContents()
            } else {
                $ws.Range($cellRef).Value = $newValue
            }
        }
    }
}

Write-Host "Applied $($updates.Keys.Count) sheet updates."